# Insert two new rows (a new weekly record) above the current row 338,
# shifting the existing data block (rows 338-473) down to rows 340-475,
# and populate the two new rows with the new week's Primera/Segunda data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 338 - everything from old row 338 onward shifts down by 2.
$ws.Rows.Item(338).Resize(2).Insert()

# --- New row 338: "Primera" quality record for the new date ---
$ws.Range("A338").Value = 11
$ws.Range("B338").Value = "Vega Monumental Concepción"
$ws.Range("C338").Value = "Bíobío"
$ws.Range("D338").Value = 45009
$ws.Range("E338").Value = 8
$ws.Range("F338").Value = 100114014
$ws.Range("G338").Value = "Betarraga"
$ws.Range("H338").Value = "Sin especificar"
$ws.Range("I338").Value = "Primera"
$ws.Range("J338").Value = 500
$ws.Range("K338").Value = 700
$ws.Range("L338").Value = 800
$ws.Range("M338").Value = 740
$ws.Range("N338").Value = "$/paquete 5 unidades"
$ws.Range("O338").Value = "Región Metropolitana"
$ws.Range("P338").Value = 148
$ws.Range("Q338").Value = 5
$ws.Range("R338").Value = "Hortaliza"

# --- New row 339: "Segunda" quality record for the new date ---
$ws.Range("A339").Value = 11
$ws.Range("B339").Value = "Vega Monumental Concepción"
$ws.Range("C339").Value = "Bíobío"
$ws.Range("D339").Value = 45009
$ws.Range("E339").Value = 8
$ws.Range("F339").Value = 100114014
$ws.Range("G339").Value = "Betarraga"
$ws.Range("H339").Value = "Sin especificar"
$ws.Range("I339").Value = "Segunda"
$ws.Range("J339").Value = 300
$ws.Range("K339").Value = 600
$ws.Range("L339").Value = 600
$ws.Range("M339").Value = 600
$ws.Range("N339").Value = "$/paquete 5 unidades"
$ws.Range("O339").Value = "Región Metropolitana"
$ws.Range("P339").Value = 120
$ws.Range("Q339").Value = 5
$ws.Range("R339").Value = "Hortaliza"
